$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns B and C
$ws.Range("B1").Value = "kode_kelas"
$ws.Range("C1").Value = "kode_kamar"

# Replace numeric placeholder codes in B/C with the real class/room codes
$ws.Range("B2").Value = "WEZNC20300"
$ws.Range("C2").Value = "QKUZQ10444"

$ws.Range("B3").Value = "WEZNC20300"
$ws.Range("C3").Value = "QKUZQ10444"

$ws.Range("B4").Value = "WEZNC20300"
$ws.Range("C4").Value = "ULKLN41341"

# Set explicit column widths for the newly-populated B and C columns
# (ColumnWidth values below are the character-width inputs that this
# engine's pixel-quantized width model stores as 10.875 / 13.5 in OOXML)
$ws.Columns.Item(2).ColumnWidth = 10.041666666666666
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666

# Update the active selection to reflect the edited cell
$ws.Range("B4").Select()
